$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns that differ across rows 2-4 and must be cyclically rotated:
# new row2 = old row4, new row3 = old row2, new row4 = old row3
$cols = @("A", "B", "D", "E", "F", "G", "H", "M", "Q", "R")

# Capture original values for rows 2, 3, 4 (use Value2 - bare Value getter
# does not resolve correctly through this COM-interop shim)
$orig2 = @{}
$orig3 = @{}
$orig4 = @{}
foreach ($col in $cols) {
    $orig2[$col] = $ws.Range("${col}2").Value2
    $orig3[$col] = $ws.Range("${col}3").Value2
    $orig4[$col] = $ws.Range("${col}4").Value2
}

# Apply rotation: row2 <- row4, row3 <- row2, row4 <- row3
foreach ($col in $cols) {
    $ws.Range("${col}2").Value = $orig4[$col]
    $ws.Range("${col}3").Value = $orig2[$col]
    $ws.Range("${col}4").Value = $orig3[$col]
}
